# Split the run containing "#mr_rec_needs___88# #mr_needs_oth#" so that a
# new, separately-formatted "Others: " label is inserted between the two
# placeholders, producing three runs (all sharing the original formatting).

$d = $word.ActiveDocument

# Locate the space right after "#mr_rec_needs___88#" and collapse the found
# range to its end (i.e. the insertion point right before "#mr_needs_oth#").
$found = $d.Content
$found.Find.Execute("#mr_rec_needs___88# ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found.Collapse(0)

# Insert the new label text at that point.
$found.InsertAfter("Others: ")

# Force Word to break the newly inserted text into its own run (distinct
# from the surrounding placeholder text) by toggling a direct character
# format on/off; since Bold is already "on" for this text, restoring it to
# "on" leaves the formatting unchanged while still causing the run split.
$found.Bold = 0
$found.Bold = 1
